$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2061.3333
$ws.Range("I12").Value = 1422.5
$ws.Range("J12").Value = 3339
$ws.Range("K12").Value = 1422.5
$ws.Range("L12").Value = 3339
$ws.Range("M12").Value = -1252.5
$ws.Range("N12").Value = -3679

$ws.Range("H40").Value = 1543.2174
$ws.Range("I40").Value = 1537.8096
$ws.Range("J40").Value = 1600
$ws.Range("K40").Value = 1537.8096
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -1362.8096
$ws.Range("N40").Value = -1950

$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138

$ws.Range("H112").Value = 2140.6667
$ws.Range("I112").Value = 954
$ws.Range("J112").Value = 2988.2856
$ws.Range("K112").Value = 2862
$ws.Range("L112").Value = 8964.856800000001
$ws.Range("M112").Value = -1754
$ws.Range("N112").Value = -11180.8568

$ws.Range("H132").Value = 2264.875
$ws.Range("I132").Value = 2082.5334
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6247.600199999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3717.600199999999
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1989.3334
$ws.Range("I8").Value = 1989.3334
$ws.Range("K8").Value = 1989.3334
$ws.Range("M8").Value = -1845.3334

$ws.Range("I22").Value = 420
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 420
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -121
$ws.Range("N22").ClearContents()

$ws.Range("H45").Value = 1744.4
$ws.Range("I45").Value = 1181.875
$ws.Range("K45").Value = 1181.875
$ws.Range("M45").Value = -804.875

$ws.Range("H61").Value = 7889.143
$ws.Range("I61").Value = 408
$ws.Range("K61").Value = 408
$ws.Range("M61").Value = -196

$ws.Range("H132").Value = 2171.5
$ws.Range("I132").Value = 2191
$ws.Range("J132").Value = 2165
$ws.Range("K132").Value = 6573
$ws.Range("L132").Value = 6495
$ws.Range("M132").Value = -4043
$ws.Range("N132").Value = -11555

$ws.Range("H136").Value = 7889.143
$ws.Range("I136").Value = 408
$ws.Range("K136").Value = 1224
$ws.Range("M136").Value = 1326

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 8
$ws.Range("I16").Value = 8
$ws.Range("K16").Value = 8
$ws.Range("M16").Value = 162

$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H134").Value = 3503.4285
$ws.Range("I134").Value = 3503.4285
$ws.Range("K134").Value = 10510.2855
$ws.Range("M134").Value = -7975.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 360
$ws.Range("I10").Value = 360
$ws.Range("K10").Value = 360
$ws.Range("M10").Value = -221

$ws.Range("H31").Value = 873.375
$ws.Range("I31").Value = 719.9286
$ws.Range("K31").Value = 719.9286
$ws.Range("M31").Value = -424.9286

$ws.Range("H34").Value = 873.375
$ws.Range("I34").Value = 719.9286
$ws.Range("K34").Value = 719.9286
$ws.Range("M34").Value = -517.9286

$ws.Range("H58").Value = 4941.75
$ws.Range("I58").Value = 2106.5
$ws.Range("J58").Value = 7777
$ws.Range("K58").Value = 2106.5
$ws.Range("L58").Value = 7777
$ws.Range("M58").Value = -1903.5
$ws.Range("N58").Value = -8183

$ws.Range("H132").Value = 7833.625
$ws.Range("I132").Value = 7833.625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 23500.875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -20970.875
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2283.25
$ws.Range("I134").Value = 2283.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6849.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4314.75
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 4941.75
$ws.Range("I136").Value = 2106.5
$ws.Range("J136").Value = 7777
$ws.Range("K136").Value = 6319.5
$ws.Range("L136").Value = 23331
$ws.Range("M136").Value = -3769.5
$ws.Range("N136").Value = -28431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9965774
$ws.Range("I4").Value = 10911477
$ws.Range("K4").Value = 32734431
$ws.Range("M4").Value = -32734319

$ws.Range("H13").Value = 5.5
$ws.Range("I13").Value = 5.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 16.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 151.5
$ws.Range("N13").ClearContents()

$ws.Range("H107").Value = 1542.375
$ws.Range("J107").Value = 1996.5
$ws.Range("L107").Value = 5989.5
$ws.Range("N107").Value = -9829.5

$ws.Range("H136").Value = 16250
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -9900
$ws.Range("N136").Value = -70200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 273.85715
$ws.Range("I2").Value = 273.4
$ws.Range("K2").Value = 273.4
$ws.Range("M2").Value = -160.4

$ws.Range("H7").Value = 1000000
$ws.Range("J7").Value = 1000000
$ws.Range("L7").Value = 1000000
$ws.Range("N7").Value = -1000224

$ws.Range("H8").Value = 1000000
$ws.Range("J8").Value = 1000000
$ws.Range("L8").Value = 1000000
$ws.Range("N8").Value = -1000278

$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730

$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064

$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 671

$ws.Range("H122").Value = 2162
$ws.Range("I122").Value = 2175.875
$ws.Range("K122").Value = 6527.625
$ws.Range("M122").Value = -4077.625

$ws.Range("H132").Value = 7144.2
$ws.Range("I132").Value = 8325.666999999999
$ws.Range("J132").Value = 5372
$ws.Range("K132").Value = 24977.001
$ws.Range("L132").Value = 16116
$ws.Range("M132").Value = -22447.001
$ws.Range("N132").Value = -21176

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6814.75
$ws.Range("J7").Value = 7899.5
$ws.Range("L7").Value = 7899.5
$ws.Range("N7").Value = -8123.5

$ws.Range("H122").Value = 5333.2954
$ws.Range("I122").Value = 4012.125
$ws.Range("J122").Value = 6918.7
$ws.Range("K122").Value = 12036.375
$ws.Range("L122").Value = 20756.1
$ws.Range("M122").Value = -9586.375
$ws.Range("N122").Value = -25656.1

$ws.Range("H126").Value = 6814.75
$ws.Range("J126").Value = 7899.5
$ws.Range("L126").Value = 23698.5
$ws.Range("N126").Value = -28638.5

$ws.Range("H132").Value = 3439.7778
$ws.Range("I132").Value = 3368.818
$ws.Range("K132").Value = 10106.454
$ws.Range("M132").Value = -7576.454000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H122").Value = 3021.5757
$ws.Range("J122").Value = 3764.2856
$ws.Range("L122").Value = 11292.8568
$ws.Range("N122").Value = -16192.8568
